$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 4447751.5  # H76: 3177811.5 -> 4447751.5
$ws.Cells.Item(76, 9).Value = 5558690.5  # I76: 3834514 -> 5558690.5
$ws.Cells.Item(76, 10).Value = 3996  # J76: 3750 -> 3996
$ws.Cells.Item(76, 11).Value = 5558690.5  # K76: 3834514 -> 5558690.5
$ws.Cells.Item(76, 12).Value = 3996  # L76: 3750 -> 3996
$ws.Cells.Item(76, 13).Value = -5558375.5  # M76: -3834199 -> -5558375.5
$ws.Cells.Item(76, 14).Value = -4626  # N76: -4380 -> -4626
$ws.Cells.Item(79, 8).Value = 4447751.5  # H79: 3177811.5 -> 4447751.5
$ws.Cells.Item(79, 9).Value = 5558690.5  # I79: 3834514 -> 5558690.5
$ws.Cells.Item(79, 10).Value = 3996  # J79: 3750 -> 3996
$ws.Cells.Item(79, 11).Value = 5558690.5  # K79: 3834514 -> 5558690.5
$ws.Cells.Item(79, 12).Value = 3996  # L79: 3750 -> 3996
$ws.Cells.Item(79, 13).Value = -5557598.5  # M79: -3833422 -> -5557598.5
$ws.Cells.Item(79, 14).Value = -6180  # N79: -5934 -> -6180
$ws.Cells.Item(98, 8).Value = 509702.5  # H98: 623167.1 -> 509702.5
$ws.Cells.Item(98, 9).Value = 622170.2  # I98: 659442.5 -> 622170.2
$ws.Cells.Item(98, 10).Value = 3598  # J98: 6486 -> 3598
$ws.Cells.Item(98, 11).Value = 622170.2  # K98: 659442.5 -> 622170.2
$ws.Cells.Item(98, 12).Value = 3598  # L98: 6486 -> 3598
$ws.Cells.Item(98, 13).Value = -620672.2  # M98: -657944.5 -> -620672.2
$ws.Cells.Item(98, 14).Value = -6594  # N98: -9482 -> -6594
$ws.Cells.Item(107, 8).Value = 617714.6  # H107: 1389739 -> 617714.6
$ws.Cells.Item(107, 9).Value = 793969.0600000001  # I107: 11111111 -> 793969.0600000001
$ws.Cells.Item(107, 10).Value = 824  # J107: 971.5714 -> 824
$ws.Cells.Item(107, 11).Value = 793969.0600000001  # K107: 11111111 -> 793969.0600000001
$ws.Cells.Item(107, 12).Value = 824  # L107: 971.5714 -> 824
$ws.Cells.Item(107, 13).Value = -792049.0600000001  # M107: -11109191 -> -792049.0600000001
$ws.Cells.Item(107, 14).Value = -4664  # N107: -4811.5714 -> -4664
$ws.Cells.Item(112, 8).Value = 11859007  # H112: 10910371 -> 11859007
$ws.Cells.Item(112, 10).Value = 12398007  # J112: 11859012 -> 12398007
$ws.Cells.Item(112, 12).Value = 37194021  # L112: 35577036 -> 37194021
$ws.Cells.Item(112, 14).Value = -37196237  # N112: -35579252 -> -37196237
$ws.Cells.Item(122, 8).Value = 509702.5  # H122: 623167.1 -> 509702.5
$ws.Cells.Item(122, 9).Value = 622170.2  # I122: 659442.5 -> 622170.2
$ws.Cells.Item(122, 10).Value = 3598  # J122: 6486 -> 3598
$ws.Cells.Item(122, 11).Value = 1866510.6  # K122: 1978327.5 -> 1866510.6
$ws.Cells.Item(122, 12).Value = 10794  # L122: 19458 -> 10794
$ws.Cells.Item(122, 13).Value = -1864060.6  # M122: -1975877.5 -> -1864060.6
$ws.Cells.Item(122, 14).Value = -15694  # N122: -24358 -> -15694
$ws.Cells.Item(137, 8).Value = 1063.8788  # H137: 1226.48 -> 1063.8788
$ws.Cells.Item(137, 9).Value = 1004.5417  # I137: 1164.5555 -> 1004.5417
$ws.Cells.Item(137, 10).Value = 1222.1111  # J137: 1385.7142 -> 1222.1111
$ws.Cells.Item(137, 11).Value = 3013.6251  # K137: 3493.6665 -> 3013.6251
$ws.Cells.Item(137, 12).Value = 3666.3333  # L137: 4157.142599999999 -> 3666.3333
$ws.Cells.Item(137, 13).Value = -463.6251000000002  # M137: -943.6664999999998 -> -463.6251000000002
$ws.Cells.Item(137, 14).Value = -8766.3333  # N137: -9257.142599999999 -> -8766.3333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(132, 8).Value = 49444  # H132: 45520 -> 49444
$ws.Cells.Item(132, 10).Value = 49444  # J132: 45520 -> 49444
$ws.Cells.Item(132, 12).Value = 49444  # L132: 45520 -> 49444
$ws.Cells.Item(132, 14).Value = -59564  # N132: -55640 -> -59564

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 500000500  # H4: 333334000 -> 500000500
$ws.Cells.Item(4, 10).Value = 500000500  # J4: 333334000 -> 500000500
$ws.Cells.Item(4, 12).Value = 500000500  # L4: 333334000 -> 500000500
$ws.Cells.Item(4, 14).Value = -500000724  # N4: -333334224 -> -500000724
$ws.Cells.Item(31, 8).Value = 5932.5757  # H31: 5763.75 -> 5932.5757
$ws.Cells.Item(31, 9).Value = 2163.353  # I31: 2136.375 -> 2163.353
$ws.Cells.Item(31, 10).Value = 9937.375  # J31: 8665.65 -> 9937.375
$ws.Cells.Item(31, 11).Value = 2163.353  # K31: 2136.375 -> 2163.353
$ws.Cells.Item(31, 12).Value = 9937.375  # L31: 8665.65 -> 9937.375
$ws.Cells.Item(31, 13).Value = -1868.353  # M31: -1841.375 -> -1868.353
$ws.Cells.Item(31, 14).Value = -10527.375  # N31: -9255.65 -> -10527.375
$ws.Cells.Item(34, 8).Value = 5932.5757  # H34: 5763.75 -> 5932.5757
$ws.Cells.Item(34, 9).Value = 2163.353  # I34: 2136.375 -> 2163.353
$ws.Cells.Item(34, 10).Value = 9937.375  # J34: 8665.65 -> 9937.375
$ws.Cells.Item(34, 11).Value = 2163.353  # K34: 2136.375 -> 2163.353
$ws.Cells.Item(34, 12).Value = 9937.375  # L34: 8665.65 -> 9937.375
$ws.Cells.Item(34, 13).Value = -1961.353  # M34: -1934.375 -> -1961.353
$ws.Cells.Item(34, 14).Value = -10341.375  # N34: -9069.65 -> -10341.375
$ws.Cells.Item(115, 8).Value = 23775.572  # H115: 25615.8 -> 23775.572
$ws.Cells.Item(115, 10).Value = 23775.572  # J115: 25615.8 -> 23775.572
$ws.Cells.Item(115, 12).Value = 23775.572  # L115: 25615.8 -> 23775.572
$ws.Cells.Item(115, 14).Value = -26125.572  # N115: -27965.8 -> -26125.572
$ws.Cells.Item(120, 8).Value = 33342  # H120: 39999.332 -> 33342
$ws.Cells.Item(120, 10).Value = 33342  # J120: 39999.332 -> 33342
$ws.Cells.Item(120, 12).Value = 33342  # L120: 39999.332 -> 33342
$ws.Cells.Item(120, 14).Value = -40600  # N120: -47257.332 -> -40600
$ws.Cells.Item(121, 8).Value = 37666.668  # H121: 45000 -> 37666.668
$ws.Cells.Item(121, 9).Value = 30000  # I121: 0 -> 30000
$ws.Cells.Item(121, 10).Value = 41500  # J121: 45000 -> 41500
$ws.Cells.Item(121, 11).Value = 30000  # K121: 0 -> 30000
$ws.Cells.Item(121, 12).Value = 41500  # L121: 45000 -> 41500
$ws.Cells.Item(121, 13).Value = -28690  # M121: None -> -28690
$ws.Cells.Item(121, 14).Value = -44120  # N121: -47620 -> -44120
$ws.Cells.Item(132, 8).Value = 3710.5334  # H132: 4030.8333 -> 3710.5334
$ws.Cells.Item(132, 9).Value = 3524.5217  # I132: 3867.7896 -> 3524.5217
$ws.Cells.Item(132, 10).Value = 4321.7144  # J132: 4650.4 -> 4321.7144
$ws.Cells.Item(132, 11).Value = 10573.5651  # K132: 11603.3688 -> 10573.5651
$ws.Cells.Item(132, 12).Value = 12965.1432  # L132: 13951.2 -> 12965.1432
$ws.Cells.Item(132, 13).Value = -8043.5651  # M132: -9073.3688 -> -8043.5651
$ws.Cells.Item(132, 14).Value = -18025.1432  # N132: -19011.2 -> -18025.1432
$ws.Cells.Item(134, 8).Value = 2385.9285  # H134: 2843.7026 -> 2385.9285
$ws.Cells.Item(134, 9).Value = 1045.1538  # I134: 1477.091 -> 1045.1538
$ws.Cells.Item(134, 10).Value = 4564.6875  # J134: 4848.067 -> 4564.6875
$ws.Cells.Item(134, 11).Value = 3135.4614  # K134: 4431.272999999999 -> 3135.4614
$ws.Cells.Item(134, 12).Value = 13694.0625  # L134: 14544.201 -> 13694.0625
$ws.Cells.Item(134, 13).Value = -600.4614000000001  # M134: -1896.272999999999 -> -600.4614000000001
$ws.Cells.Item(134, 14).Value = -18764.0625  # N134: -19614.201 -> -18764.0625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 10450.862  # H4: 9482.968999999999 -> 10450.862
$ws.Cells.Item(4, 9).Value = 94.21429000000001  # I4: 99.92308 -> 94.21429000000001
$ws.Cells.Item(4, 10).Value = 20117.066  # J4: 15902.947 -> 20117.066
$ws.Cells.Item(4, 11).Value = 282.64287  # K4: 299.76924 -> 282.64287
$ws.Cells.Item(4, 12).Value = 60351.198  # L4: 47708.841 -> 60351.198
$ws.Cells.Item(4, 13).Value = -170.64287  # M4: -187.76924 -> -170.64287
$ws.Cells.Item(4, 14).Value = -60575.198  # N4: -47932.841 -> -60575.198
$ws.Cells.Item(107, 8).Value = 394.7857  # H107: 360.3243 -> 394.7857
$ws.Cells.Item(107, 9).Value = 366.8095  # I107: 356.4091 -> 366.8095
$ws.Cells.Item(107, 10).Value = 478.7143  # J107: 366.06668 -> 478.7143
$ws.Cells.Item(107, 11).Value = 1100.4285  # K107: 1069.2273 -> 1100.4285
$ws.Cells.Item(107, 12).Value = 1436.1429  # L107: 1098.20004 -> 1436.1429
$ws.Cells.Item(107, 13).Value = 819.5715  # M107: 850.7727 -> 819.5715
$ws.Cells.Item(107, 14).Value = -5276.1429  # N107: -4938.20004 -> -5276.1429
$ws.Cells.Item(113, 8).Value = 826.1875  # H113: 945.1429000000001 -> 826.1875
$ws.Cells.Item(113, 9).Value = 699.1667  # I113: 0 -> 699.1667
$ws.Cells.Item(113, 10).Value = 902.4  # J113: 945.1429000000001 -> 902.4
$ws.Cells.Item(113, 11).Value = 2097.5001  # K113: 0 -> 2097.5001
$ws.Cells.Item(113, 12).Value = 2707.2  # L113: 2835.4287 -> 2707.2
$ws.Cells.Item(113, 13).Value = 72.4998999999998  # M113: None -> 72.4998999999998
$ws.Cells.Item(113, 14).Value = -7047.2  # N113: -7175.4287 -> -7047.2
$ws.Cells.Item(117, 8).Value = 713.3  # H117: 1000 -> 713.3
$ws.Cells.Item(117, 9).Value = 0  # I117: 200 -> 0
$ws.Cells.Item(117, 10).Value = 713.3  # J117: 1400 -> 713.3
$ws.Cells.Item(117, 11).Value = 0  # K117: 600 -> 0
$ws.Cells.Item(117, 12).Value = 2139.9  # L117: 4200 -> 2139.9
$ws.Cells.Item(117, 13).ClearContents()  # M117: 2842 -> (removed)
$ws.Cells.Item(117, 14).Value = -9023.9  # N117: -11084 -> -9023.9
$ws.Cells.Item(129, 8).Value = 957.1667  # H129: 494.08334 -> 957.1667
$ws.Cells.Item(129, 9).Value = 725  # I129: 342.9 -> 725
$ws.Cells.Item(129, 10).Value = 1073.25  # J129: 1250 -> 1073.25
$ws.Cells.Item(129, 11).Value = 2175  # K129: 1028.7 -> 2175
$ws.Cells.Item(129, 12).Value = 3219.75  # L129: 3750 -> 3219.75
$ws.Cells.Item(129, 13).Value = 2825  # M129: 3971.3 -> 2825
$ws.Cells.Item(129, 14).Value = -13219.75  # N129: -13750 -> -13219.75
$ws.Cells.Item(131, 8).Value = 6945864  # H131: 7408834.5 -> 6945864
$ws.Cells.Item(131, 9).Value = 386.33334  # I131: 369.7143 -> 386.33334
$ws.Cells.Item(131, 10).Value = 7938075  # J131: 8773552 -> 7938075
$ws.Cells.Item(131, 11).Value = 1159.00002  # K131: 1109.1429 -> 1159.00002
$ws.Cells.Item(131, 12).Value = 23814225  # L131: 26320656 -> 23814225
$ws.Cells.Item(131, 13).Value = 3880.99998  # M131: 3930.8571 -> 3880.99998
$ws.Cells.Item(131, 14).Value = -23824305  # N131: -26330736 -> -23824305

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 35717130  # H80: 29414570 -> 35717130
$ws.Cells.Item(80, 9).Value = 2999.0908  # I80: 2929.6155 -> 2999.0908
$ws.Cells.Item(80, 10).Value = 166668930  # J80: 125002400 -> 166668930
$ws.Cells.Item(80, 11).Value = 2999.0908  # K80: 2929.6155 -> 2999.0908
$ws.Cells.Item(80, 12).Value = 166668930  # L80: 125002400 -> 166668930
$ws.Cells.Item(80, 13).Value = -2001.0908  # M80: -1931.6155 -> -2001.0908
$ws.Cells.Item(80, 14).Value = -166670926  # N80: -125004396 -> -166670926
$ws.Cells.Item(83, 8).Value = 35717130  # H83: 29414570 -> 35717130
$ws.Cells.Item(83, 9).Value = 2999.0908  # I83: 2929.6155 -> 2999.0908
$ws.Cells.Item(83, 10).Value = 166668930  # J83: 125002400 -> 166668930
$ws.Cells.Item(83, 11).Value = 14995.454  # K83: 14648.0775 -> 14995.454
$ws.Cells.Item(83, 12).Value = 833344650  # L83: 625012000 -> 833344650
$ws.Cells.Item(83, 13).Value = -10003.454  # M83: -9656.077499999999 -> -10003.454
$ws.Cells.Item(83, 14).Value = -833354634  # N83: -625021984 -> -833354634
$ws.Cells.Item(113, 8).Value = 1199.3334  # H113: 1265.9333 -> 1199.3334
$ws.Cells.Item(113, 9).Value = 1152.5333  # I113: 1220.6428 -> 1152.5333
$ws.Cells.Item(113, 10).Value = 1433.3334  # J113: 1900 -> 1433.3334
$ws.Cells.Item(113, 11).Value = 1152.5333  # K113: 1220.6428 -> 1152.5333
$ws.Cells.Item(113, 12).Value = 1433.3334  # L113: 1900 -> 1433.3334
$ws.Cells.Item(113, 13).Value = 1017.4667  # M113: 949.3571999999999 -> 1017.4667
$ws.Cells.Item(113, 14).Value = -5773.3334  # N113: -6240 -> -5773.3334
$ws.Cells.Item(126, 8).Value = 2667.5  # H126: 2559.3142 -> 2667.5
$ws.Cells.Item(126, 9).Value = 2258.3333  # I126: 2092.7058 -> 2258.3333
$ws.Cells.Item(126, 10).Value = 2842.8572  # J126: 3000 -> 2842.8572
$ws.Cells.Item(126, 11).Value = 6774.999899999999  # K126: 6278.117400000001 -> 6774.999899999999
$ws.Cells.Item(126, 12).Value = 8528.571599999999  # L126: 9000 -> 8528.571599999999
$ws.Cells.Item(126, 13).Value = -4304.999899999999  # M126: -3808.117400000001 -> -4304.999899999999
$ws.Cells.Item(126, 14).Value = -13468.5716  # N126: -13940 -> -13468.5716

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 31000  # H109: 32000 -> 31000
$ws.Cells.Item(109, 10).Value = 31000  # J109: 32000 -> 31000
$ws.Cells.Item(109, 12).Value = 31000  # L109: 32000 -> 31000
$ws.Cells.Item(109, 14).Value = -33774  # N109: -34774 -> -33774
$ws.Cells.Item(113, 8).Value = 426.58334  # H113: 260.90475 -> 426.58334
$ws.Cells.Item(113, 9).Value = 374.45456  # I113: 233.10527 -> 374.45456
$ws.Cells.Item(113, 10).Value = 1000  # J113: 525 -> 1000
$ws.Cells.Item(113, 11).Value = 1123.36368  # K113: 699.3158099999999 -> 1123.36368
$ws.Cells.Item(113, 12).Value = 3000  # L113: 1575 -> 3000
$ws.Cells.Item(113, 13).Value = 1046.63632  # M113: 1470.68419 -> 1046.63632
$ws.Cells.Item(113, 14).Value = -7340  # N113: -5915 -> -7340
$ws.Cells.Item(122, 8).Value = 112600.22  # H122: 72512.28999999999 -> 112600.22
$ws.Cells.Item(122, 9).Value = 126425.25  # I122: 167418.67 -> 126425.25
$ws.Cells.Item(122, 10).Value = 2000  # J122: 1332.5 -> 2000
$ws.Cells.Item(122, 11).Value = 379275.75  # K122: 502256.01 -> 379275.75
$ws.Cells.Item(122, 12).Value = 6000  # L122: 3997.5 -> 6000
$ws.Cells.Item(122, 13).Value = -376825.75  # M122: -499806.01 -> -376825.75
$ws.Cells.Item(122, 14).Value = -10900  # N122: -8897.5 -> -10900
$ws.Cells.Item(132, 8).Value = 2055.5898  # H132: 2345.4243 -> 2055.5898
$ws.Cells.Item(132, 9).Value = 1767.6562  # I132: 2096.1482 -> 1767.6562
$ws.Cells.Item(132, 10).Value = 3371.8572  # J132: 3467.1667 -> 3371.8572
$ws.Cells.Item(132, 11).Value = 5302.9686  # K132: 6288.444600000001 -> 5302.9686
$ws.Cells.Item(132, 12).Value = 10115.5716  # L132: 10401.5001 -> 10115.5716
$ws.Cells.Item(132, 13).Value = -2772.9686  # M132: -3758.444600000001 -> -2772.9686
$ws.Cells.Item(132, 14).Value = -15175.5716  # N132: -15461.5001 -> -15175.5716
$ws.Cells.Item(136, 8).Value = 1762.4073  # H136: 1949.5416 -> 1762.4073
$ws.Cells.Item(136, 9).Value = 815.8333  # I136: 925.93335 -> 815.8333
$ws.Cells.Item(136, 11).Value = 2447.4999  # K136: 2777.80005 -> 2447.4999
$ws.Cells.Item(136, 13).Value = 102.5001000000002  # M136: -227.8000499999998 -> 102.5001000000002

Write-Output "Applied 196 cell updates across ALC, BSM, CRP, CUL, GSM, WVR sheets"
